# renamed repo, fixed output folder path
# even_MAG-GUT17212.fa was removed from the dataset; delete its row (row 4)
# and shift the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Delete()
